$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Swap match data (columns F:V) between rows 14 and 15
# -----------------------------------------------------------------
$row14 = $ws.Range("F14:V14").Value()
$row15 = $ws.Range("F15:V15").Value()
$ws.Range("F14:V14").Value = $row15
$ws.Range("F15:V15").Value = $row14

# -----------------------------------------------------------------
# 2) Swap match data (columns F:V) between rows 108 and 109
# -----------------------------------------------------------------
$row108 = $ws.Range("F108:V108").Value()
$row109 = $ws.Range("F109:V109").Value()
$ws.Range("F108:V108").Value = $row109
$ws.Range("F109:V109").Value = $row108

# -----------------------------------------------------------------
# 3) Append three new match rows (119, 120, 121) after the last
#    existing row (118). Copy row 118 first so the new rows pick
#    up the exact same cell styles (bold/border on col A, date
#    format on col E), then overwrite the values cell by cell.
# -----------------------------------------------------------------
$ws.Range("A118:V118").Copy($ws.Range("A119:V119"))
$ws.Range("A118:V118").Copy($ws.Range("A120:V120"))
$ws.Range("A118:V118").Copy($ws.Range("A121:V121"))

$ws.Cells.Item(119, 1).Value = 118
$ws.Cells.Item(119, 2).Value = "czech-republic"
$ws.Cells.Item(119, 3).Value = "fnl"
$ws.Cells.Item(119, 4).Value = "2023-2024"
$ws.Cells.Item(119, 5).Value = 45235.42708333334
$ws.Cells.Item(119, 6).Value = "Sigma Olomouc B"
$ws.Cells.Item(119, 7).Value = 2
$ws.Cells.Item(119, 8).Value = "Chrudim"
$ws.Cells.Item(119, 9).Value = 4
$ws.Cells.Item(119, 10).Value = 2.11
$ws.Cells.Item(119, 11).Value = "05/11/2023 02:23"
$ws.Cells.Item(119, 12).Value = 2.11
$ws.Cells.Item(119, 13).Value = "05/11/2023 02:23"
$ws.Cells.Item(119, 14).Value = 3.24
$ws.Cells.Item(119, 15).Value = "05/11/2023 08:20"
$ws.Cells.Item(119, 16).Value = 3.24
$ws.Cells.Item(119, 17).Value = "05/11/2023 08:20"
$ws.Cells.Item(119, 18).Value = 3.57
$ws.Cells.Item(119, 19).Value = "05/11/2023 02:23"
$ws.Cells.Item(119, 20).Value = 3.57
$ws.Cells.Item(119, 21).Value = "05/11/2023 02:23"
$ws.Cells.Item(119, 22).Value = "https://www.betexplorer.com/football/czech-republic/fnl/sigma-olomouc-chrudim/fVu5dKI6/"

$ws.Cells.Item(120, 1).Value = 119
$ws.Cells.Item(120, 2).Value = "czech-republic"
$ws.Cells.Item(120, 3).Value = "fnl"
$ws.Cells.Item(120, 4).Value = "2023-2024"
$ws.Cells.Item(120, 5).Value = 45235.42708333334
$ws.Cells.Item(120, 6).Value = "Zizkov"
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(120, 8).Value = "Jihlava"
$ws.Cells.Item(120, 9).Value = 2
$ws.Cells.Item(120, 10).Value = 1.81
$ws.Cells.Item(120, 11).Value = "05/11/2023 08:52"
$ws.Cells.Item(120, 12).Value = 1.81
$ws.Cells.Item(120, 13).Value = "05/11/2023 08:52"
$ws.Cells.Item(120, 14).Value = 3.78
$ws.Cells.Item(120, 15).Value = "05/11/2023 08:52"
$ws.Cells.Item(120, 16).Value = 3.78
$ws.Cells.Item(120, 17).Value = "05/11/2023 08:52"
$ws.Cells.Item(120, 18).Value = 4.18
$ws.Cells.Item(120, 19).Value = "05/11/2023 08:52"
$ws.Cells.Item(120, 20).Value = 4.18
$ws.Cells.Item(120, 21).Value = "05/11/2023 08:52"
$ws.Cells.Item(120, 22).Value = "https://www.betexplorer.com/football/czech-republic/fnl/zizkov-jihlava/S2Me02Js/"

$ws.Cells.Item(121, 1).Value = 120
$ws.Cells.Item(121, 2).Value = "czech-republic"
$ws.Cells.Item(121, 3).Value = "fnl"
$ws.Cells.Item(121, 4).Value = "2023-2024"
$ws.Cells.Item(121, 5).Value = 45235.58333333334
$ws.Cells.Item(121, 6).Value = "Dukla Prague"
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = "Taborsko"
$ws.Cells.Item(121, 9).Value = 0
$ws.Cells.Item(121, 10).Value = 2.08
$ws.Cells.Item(121, 11).Value = "03/11/2023 00:12"
$ws.Cells.Item(121, 12).Value = 1.75
$ws.Cells.Item(121, 13).Value = "05/11/2023 13:54"
$ws.Cells.Item(121, 14).Value = 3.4
$ws.Cells.Item(121, 15).Value = "03/11/2023 00:12"
$ws.Cells.Item(121, 16).Value = 3.88
$ws.Cells.Item(121, 17).Value = "05/11/2023 13:54"
$ws.Cells.Item(121, 18).Value = 3.06
$ws.Cells.Item(121, 19).Value = "03/11/2023 00:12"
$ws.Cells.Item(121, 20).Value = 4.38
$ws.Cells.Item(121, 21).Value = "05/11/2023 13:54"
$ws.Cells.Item(121, 22).Value = "https://www.betexplorer.com/football/czech-republic/fnl/dukla-prague-taborsko/xfV9evYC/"

"done"
